$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5060
$ws1.Range("F13").Value = 1413
$ws1.Range("F14").Value = 3739
$ws1.Range("F16").Value = 144
$ws1.Range("F17").Value = 130
$ws1.Range("F19").Value = 2731
$ws1.Range("F20").Value = 136
$ws1.Range("F21").Value = 30
$ws1.Range("F26").Value = 8
$ws1.Range("F30").Value = 46

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 5060
$ws4.Range("F14").Value = 1413
$ws4.Range("F15").Value = 3739
$ws4.Range("F17").Value = 144
$ws4.Range("F18").Value = 130
$ws4.Range("F20").Value = 2731
$ws4.Range("F21").Value = 136
$ws4.Range("F22").Value = 30
$ws4.Range("F27").Value = 8
$ws4.Range("F31").Value = 46
